# "agregue el whatsapp en editar trabajo"
# Registros.xlsx - "Obras en general" sheet:
#  - clear N87/O87/P87 (visado de gas/salubridad/electrica) which were
#    mistakenly filled while editing a job
#  - add two new job rows (97, 98) that were registered afterwards

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Obras en general")

# --- clear the 3 visado values that were removed from row 87 ---
$ws.Cells.Item(87, 14).Value = $null   # N87
$ws.Cells.Item(87, 15).Value = $null   # O87
$ws.Cells.Item(87, 16).Value = $null   # P87

# --- row 97: new "Obra nueva" job ---
$ws.Cells.Item(97, 1).Value = "'05/06/2025"              # A - Fecha
$ws.Cells.Item(97, 2).Value = "MMO"                       # B - Profesion
$ws.Cells.Item(97, 3).Value = "Digital"                   # C - Formato
$ws.Cells.Item(97, 5).Value = "Obra nueva"                # E - Tipo de trabajo
$ws.Cells.Item(97, 6).Value = "CORDOVES MAURO IVAN"       # F - Profesional
$ws.Cells.Item(97, 7).Value = "BUSTOS KARINA VANESSA"     # G - Comitente
$ws.Cells.Item(97, 8).Value = "CALLE (47) M. HERRERA N° 4907"  # H - Ubicacion
$ws.Cells.Item(97, 9).Value = "'34987/2024"               # I - Nro expte municipal
$ws.Cells.Item(97, 10).Value = "'4072"                    # J - Nro de sistema GOP
$ws.Cells.Item(97, 11).Value = "'22090"                   # K - Nro de partida inmobiliaria
$ws.Cells.Item(97, 18).Value = "No pagado"                # R - Estado pago sellado
$ws.Cells.Item(97, 19).Value = "No pagado"                # S - Estado pago visado
$ws.Cells.Item(97, 24).Value = "\\DESKTOP-5KNILLM\Users\Usuario\Compartidos\cpim_sistema\dist\Sistema CPIM\trabajos\OBRA NUEVA\CORDOVES MAURO IVAN\BUSTOS KARINA VANESSA"  # X - Ruta de carpeta

# --- row 98: new "Registracion" job ---
$ws.Cells.Item(98, 1).Value = "'06/06/2025"               # A - Fecha
$ws.Cells.Item(98, 2).Value = "Ingeniero"                 # B - Profesion
$ws.Cells.Item(98, 3).Value = "Digital"                   # C - Formato
$ws.Cells.Item(98, 5).Value = "Registración"              # E - Tipo de trabajo
$ws.Cells.Item(98, 6).Value = "EDSBERG IVAN"              # F - Profesional
$ws.Cells.Item(98, 7).Value = "CORTES MARIA DANIELLA"     # G - Comitente
$ws.Cells.Item(98, 8).Value = "Lanusse N° 2072"           # H - Ubicacion
$ws.Cells.Item(98, 9).Value = "'15492/M/1996"             # I - Nro expte municipal
$ws.Cells.Item(98, 10).Value = "'4574"                    # J - Nro de sistema GOP
$ws.Cells.Item(98, 11).Value = "'16465"                   # K - Nro de partida inmobiliaria
$ws.Cells.Item(98, 18).Value = "No pagado"                # R - Estado pago sellado
$ws.Cells.Item(98, 19).Value = "No pagado"                # S - Estado pago visado
$ws.Cells.Item(98, 24).Value = "\\DESKTOP-5KNILLM\Users\Usuario\Compartidos\cpim_sistema\dist\Sistema CPIM\trabajos\REGISTRACION\EDSBERG IVAN\CORTES MARIA DANIELLA"  # X - Ruta de carpeta
